# B6-PowerPoint.pptx edit: 2020-06-07
#
# 1) Re-style the three data tables (slides 14, 15, 16) from the
#    "Table_0" custom table style to the built-in table style
#    {09C339AD-D5EC-45EC-A385-E137D7461BFD}.
# 2) Swap the presentation's visual theme palette to the stock
#    "Office" color scheme (best-effort: only the ThemeColorScheme
#    is reachable through the exposed PowerPoint object model).

$p = $ppt.ActivePresentation

$newTableStyle = "{09C339AD-D5EC-45EC-A385-E137D7461BFD}"
foreach ($slideIndex in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIndex)
    foreach ($shape in $slide.Shapes) {
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyle)
        }
    }
}
